# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" (sheet 1) and "全部类型" (sheet 4) worksheets to reflect the
# refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5715
$wsExhibition.Range("F3").Value = 87
$wsExhibition.Range("F4").Value = 10
$wsExhibition.Range("F7").Value = 2628
$wsExhibition.Range("F9").Value = 189
$wsExhibition.Range("F11").Value = 100
$wsExhibition.Range("F12").Value = 40
$wsExhibition.Range("F13").Value = 2462
$wsExhibition.Range("F14").Value = 508

# --- Sheet "全部类型" --------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5715
$wsAll.Range("F3").Value = 87
$wsAll.Range("F4").Value = 10
$wsAll.Range("F8").Value = 2628
$wsAll.Range("F10").Value = 189
$wsAll.Range("F13").Value = 100
$wsAll.Range("F14").Value = 40
$wsAll.Range("F15").Value = 2462
$wsAll.Range("F16").Value = 508
